# Apply updated cryptocurrency data (price + volume%) to the cryptos worksheet.
# Rows 41/42 also swap coin identity (MultiversX <-> Celestia reordered).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $value) {
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell "D2" "42.556.06"
Set-TextCell "E2" "  -1.49%  "
Set-TextCell "D3" "2.236.98"
Set-TextCell "E3" "  -1.79%  "
Set-TextCell "E4" "  +0.11%  "
Set-TextCell "D5" "115.30"
Set-TextCell "E5" "  +2.91%  "
Set-TextCell "D6" "281.92"
Set-TextCell "E6" "  +6.61%  "
Set-TextCell "E7" "  -2.18%  "
Set-TextCell "E8" "  +0.22%  "
Set-TextCell "D9" "0.611"
Set-TextCell "E9" "  +0.35%  "
Set-TextCell "D10" "46.79"
Set-TextCell "E10" "  +0.28%  "
Set-TextCell "D11" "0.0929"
Set-TextCell "E11" "  -0.89%  "
Set-TextCell "D12" "9.16"
Set-TextCell "E12" "  -0.63%  "
Set-TextCell "E13" "  -2.86%  "
Set-TextCell "D14" "15.30"
Set-TextCell "E14" "  -0.46%  "
Set-TextCell "D15" "0.882"
Set-TextCell "E15" "  +2.32%  "
Set-TextCell "D16" "2.572.89"
Set-TextCell "E16" "  -1.81%  "
Set-TextCell "D17" "2.233.00"
Set-TextCell "E17" "  -1.57%  "
Set-TextCell "D18" "42.738.04"
Set-TextCell "E18" "  -1.05%  "
Set-TextCell "E19" "  -0.79%  "
Set-TextCell "D20" "6.86"
Set-TextCell "E20" "  +1.69%  "
Set-TextCell "D21" "72.20"
Set-TextCell "E21" "  +0.03%  "
Set-TextCell "E22" "  -3.82%  "
Set-TextCell "D23" "3.09"
Set-TextCell "E23" "  +7.90%  "
Set-TextCell "D24" "231.85"
Set-TextCell "E24" "  -1.10%  "
Set-TextCell "D25" "9.33"
Set-TextCell "E25" "  -0.48%  "
Set-TextCell "D26" "12.06"
Set-TextCell "E26" "  +6.22%  "
Set-TextCell "E27" "  -1.64%  "
Set-TextCell "D28" "40.38"
Set-TextCell "E28" "  -2.72%  "
Set-TextCell "E29" "  -0.27%  "
Set-TextCell "D30" "3.28"
Set-TextCell "E30" "  -2.10%  "
Set-TextCell "D31" "173.88"
Set-TextCell "E31" "  +0.38%  "
Set-TextCell "D32" "21.15"
Set-TextCell "E32" "  -1.66%  "
Set-TextCell "D33" "0.0898"
Set-TextCell "E33" "  +0.18%  "
Set-TextCell "D34" "4.49"
Set-TextCell "E34" "  +15.62%  "
Set-TextCell "D35" "5.59"
Set-TextCell "E35" "  -1.43%  "
Set-TextCell "D37" "0.0372"
Set-TextCell "E37" "  -2.09%  "
Set-TextCell "D38" "4.65"
Set-TextCell "E38" "  -0.70%  "
Set-TextCell "E39" "  +1.71%  "
Set-TextCell "E40" "  +0.33%  "
Set-TextCell "B41" "Celestia"
Set-TextCell "C41" "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextCell "D41" "13.53"
Set-TextCell "E41" "  -5.58%  "
Set-TextCell "B42" "MultiversX"
Set-TextCell "C42" "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextCell "D42" "72.09"
Set-TextCell "E42" "  -4.58%  "
Set-TextCell "D43" "0.235"
Set-TextCell "E43" "  -0.50%  "
Set-TextCell "E44" "  +0.37%  "
Set-TextCell "D45" "1.34"
Set-TextCell "E45" "  -1.64%  "
Set-TextCell "D46" "5.55"
Set-TextCell "E46" "  -9.08%  "
Set-TextCell "E47" "  +1.04%  "
Set-TextCell "D48" "8.48"
Set-TextCell "E48" "  -0.82%  "
Set-TextCell "D49" "0.652"
Set-TextCell "E49" "  +9.79%  "
Set-TextCell "D50" "0.0986"
Set-TextCell "E50" "  -0.54%  "
Set-TextCell "D51" "100.71"
Set-TextCell "E51" "  +0.28%  "
